$d = $word.ActiveDocument

# 1. Update experience years from '21 years' to '15+ years' in the professional summary.
$d.Content.Find.Execute("21 years of experience", $true, $false, $false, $false, $false, $true, 1, $false, "15+ years of experience", 2)

# 2. Remove the EDUCATION section entirely: the "EDUCATION" Heading2 paragraph plus the
#    two Heading3 degree paragraphs that follow it (Master's and Bachelor's in Political Science).
$eduHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n", "`a") -eq "EDUCATION") {
        $eduHeading = $p
        break
    }
}

if ($eduHeading -ne $null) {
    $startPos = $eduHeading.Range.Start
    $endPos = $startPos
    $scanStart = $eduHeading.Range.End
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -ge $scanStart) {
            $text = $p.Range.Text.TrimEnd("`r", "`n", "`a")
            if ($text -eq "Master of Arts in Political Science - University of California, Berkeley" -or
                $text -eq "Bachelor of Arts in Political Science - University of California, Berkeley") {
                $endPos = $p.Range.End
            } else {
                break
            }
        }
    }
    $range = $d.Range($startPos, $endPos)
    $range.Delete()
}
